$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 5).Value = "'2025/11/03"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("D2").Copy($ws.Range("E2"))
$ws.Cells.Item(2, 5).Value = "上证"
$ws.Range("D3").Copy($ws.Range("E3"))
$ws.Cells.Item(3, 5).Value = 62.24
$ws.Range("D4").Copy($ws.Range("E4"))
$ws.Cells.Item(4, 5).Value = 3956.72
$ws.Range("D5").Copy($ws.Range("E5"))
$ws.Range("D6").Copy($ws.Range("E6"))
$ws.Cells.Item(6, 5).Value = 49.02
$ws.Range("D7").Copy($ws.Range("E7"))
$ws.Cells.Item(7, 5).Value = 5545.87
$ws.Range("D8").Copy($ws.Range("E8"))
$ws.Range("D9").Copy($ws.Range("E9"))
$ws.Cells.Item(9, 5).Value = 53.94
$ws.Range("D10").Copy($ws.Range("E10"))
$ws.Cells.Item(10, 5).Value = 4619.67
$ws.Range("D11").Copy($ws.Range("E11"))
$ws.Range("D12").Copy($ws.Range("E12"))
$ws.Cells.Item(12, 5).Value = 57.54
$ws.Range("D13").Copy($ws.Range("E13"))
$ws.Cells.Item(13, 5).Value = 7255.42
$ws.Range("D14").Copy($ws.Range("E14"))
$ws.Range("D15").Copy($ws.Range("E15"))
$ws.Cells.Item(15, 5).Value = 26.75
$ws.Range("D16").Copy($ws.Range("E16"))
$ws.Cells.Item(16, 5).Value = 2699.03
$ws.Range("D17").Copy($ws.Range("E17"))
$ws.Range("D18").Copy($ws.Range("E18"))
$ws.Cells.Item(18, 5).Value = 96.63
$ws.Range("D19").Copy($ws.Range("E19"))
$ws.Cells.Item(19, 5).Value = 6840.2
$ws.Range("D20").Copy($ws.Range("E20"))
$ws.Range("D21").Copy($ws.Range("E21"))
$ws.Cells.Item(21, 5).Value = 65.76000000000001
$ws.Range("D22").Copy($ws.Range("E22"))
$ws.Cells.Item(22, 5).Value = 83935.38
$ws.Range("D23").Copy($ws.Range("E23"))
$ws.Range("D24").Copy($ws.Range("E24"))
$ws.Cells.Item(24, 5).Value = 85.7
$ws.Range("D25").Copy($ws.Range("E25"))
$ws.Cells.Item(25, 5).Value = 19909.14
$ws.Range("D26").Copy($ws.Range("E26"))
$ws.Range("D27").Copy($ws.Range("E27"))
$ws.Cells.Item(27, 5).Value = 83.77
$ws.Range("D28").Copy($ws.Range("E28"))
$ws.Cells.Item(28, 5).Value = 39894.54
$ws.Range("D29").Copy($ws.Range("E29"))
$ws.Range("D30").Copy($ws.Range("E30"))
$ws.Cells.Item(30, 5).Value = 57.46
$ws.Range("D31").Copy($ws.Range("E31"))
$ws.Cells.Item(31, 5).Value = 5665.89
$ws.Range("D32").Copy($ws.Range("E32"))
$ws.Range("D33").Copy($ws.Range("E33"))
$ws.Cells.Item(33, 5).Value = 11.77
$ws.Range("D34").Copy($ws.Range("E34"))
$ws.Cells.Item(34, 5).Value = 33010.92
$ws.Range("D35").Copy($ws.Range("E35"))
$ws.Range("D36").Copy($ws.Range("E36"))
$ws.Cells.Item(36, 5).Value = 28.98
$ws.Range("D37").Copy($ws.Range("E37"))
$ws.Cells.Item(37, 5).Value = 3321.04
$ws.Range("D38").Copy($ws.Range("E38"))
$ws.Range("D39").Copy($ws.Range("E39"))
$ws.Cells.Item(39, 5).Value = 47.3
$ws.Range("D40").Copy($ws.Range("E40"))
$ws.Cells.Item(40, 5).Value = 3143.95
$ws.Range("D41").Copy($ws.Range("E41"))
$ws.Range("D42").Copy($ws.Range("E42"))
$ws.Cells.Item(42, 5).Value = 19.05
$ws.Range("D43").Copy($ws.Range("E43"))
$ws.Cells.Item(43, 5).Value = 7354.49
$ws.Range("D44").Copy($ws.Range("E44"))
$ws.Range("D45").Copy($ws.Range("E45"))
$ws.Cells.Item(45, 5).Value = 32.5
$ws.Range("D46").Copy($ws.Range("E46"))
$ws.Cells.Item(46, 5).Value = 8889.860000000001
$ws.Range("D47").Copy($ws.Range("E47"))
$ws.Range("D48").Copy($ws.Range("E48"))
$ws.Cells.Item(48, 5).Value = 13.2
$ws.Range("D49").Copy($ws.Range("E49"))
$ws.Cells.Item(49, 5).Value = 12997.93
$ws.Range("D50").Copy($ws.Range("E50"))
$ws.Range("D51").Copy($ws.Range("E51"))
$ws.Cells.Item(51, 5).Value = 24.96
$ws.Range("D52").Copy($ws.Range("E52"))
$ws.Cells.Item(52, 5).Value = 12426.7
$ws.Range("D53").Copy($ws.Range("E53"))
$ws.Range("D54").Copy($ws.Range("E54"))
$ws.Cells.Item(54, 5).Value = 21.45
$ws.Range("D55").Copy($ws.Range("E55"))
$ws.Cells.Item(55, 5).Value = 9784.74
$ws.Range("D56").Copy($ws.Range("E56"))
$ws.Range("D57").Copy($ws.Range("E57"))
$ws.Cells.Item(57, 5).Value = 26.77
$ws.Range("D58").Copy($ws.Range("E58"))
$ws.Cells.Item(58, 5).Value = 16111.26
$ws.Range("D59").Copy($ws.Range("E59"))
$ws.Range("D60").Copy($ws.Range("E60"))
$ws.Cells.Item(60, 5).Value = 32.17
$ws.Range("D61").Copy($ws.Range("E61"))
$ws.Cells.Item(61, 5).Value = 17526.85
$ws.Range("D62").Copy($ws.Range("E62"))
$ws.Range("D63").Copy($ws.Range("E63"))
$ws.Cells.Item(63, 5).Value = 20.84
$ws.Range("D64").Copy($ws.Range("E64"))
$ws.Cells.Item(64, 5).Value = 10263.12
$ws.Range("D65").Copy($ws.Range("E65"))
$ws.Range("D66").Copy($ws.Range("E66"))
$ws.Cells.Item(66, 5).Value = 14.8
$ws.Range("D67").Copy($ws.Range("E67"))
$ws.Cells.Item(67, 5).Value = 9907.07
$ws.Range("D68").Copy($ws.Range("E68"))
$ws.Range("D69").Copy($ws.Range("E69"))
$ws.Cells.Item(69, 5).Value = 20.88
$ws.Range("D70").Copy($ws.Range("E70"))
$ws.Cells.Item(70, 5).Value = 3164.02
$ws.Range("D71").Copy($ws.Range("E71"))
$ws.Range("D72").Copy($ws.Range("E72"))
$ws.Cells.Item(72, 5).Value = 43.63
$ws.Range("D73").Copy($ws.Range("E73"))
$ws.Cells.Item(73, 5).Value = 5894.12
$ws.Range("D74").Copy($ws.Range("E74"))
$ws.Range("D75").Copy($ws.Range("E75"))
$ws.Cells.Item(75, 5).Value = 29.82
$ws.Range("D76").Copy($ws.Range("E76"))
$ws.Cells.Item(76, 5).Value = 9444.360000000001
$ws.Range("D77").Copy($ws.Range("E77"))
$ws.Range("D78").Copy($ws.Range("E78"))
$ws.Cells.Item(78, 5).Value = 18.01
$ws.Range("D79").Copy($ws.Range("E79"))
$ws.Cells.Item(79, 5).Value = 2422.48
$ws.Range("D80").Copy($ws.Range("E80"))
$ws.Range("D81").Copy($ws.Range("E81"))
$ws.Cells.Item(81, 5).Value = 56.06
$ws.Range("D82").Copy($ws.Range("E82"))
$ws.Cells.Item(82, 5).Value = 2851.5
$ws.Range("D83").Copy($ws.Range("E83"))
$ws.Range("D84").Copy($ws.Range("E84"))
$ws.Cells.Item(84, 5).Value = 58.79
$ws.Range("D85").Copy($ws.Range("E85"))
$ws.Cells.Item(85, 5).Value = 2966.57
$ws.Range("D86").Copy($ws.Range("E86"))
$ws.Range("D87").Copy($ws.Range("E87"))
$ws.Cells.Item(87, 5).Value = 52.78
$ws.Range("D88").Copy($ws.Range("E88"))
$ws.Cells.Item(88, 5).Value = 3909.78
$ws.Range("D89").Copy($ws.Range("E89"))
$ws.Range("D90").Copy($ws.Range("E90"))
$ws.Cells.Item(90, 5).Value = 47.79
$ws.Range("D91").Copy($ws.Range("E91"))
$ws.Cells.Item(91, 5).Value = 2083.67
$ws.Range("D92").Copy($ws.Range("E92"))
$ws.Range("D93").Copy($ws.Range("E93"))
$ws.Cells.Item(93, 5).Value = 28.1
$ws.Range("D94").Copy($ws.Range("E94"))
$ws.Cells.Item(94, 5).Value = 13963.26
$ws.Range("D95").Copy($ws.Range("E95"))
$ws.Range("D96").Copy($ws.Range("E96"))
$ws.Cells.Item(96, 5).Value = 86.69
$ws.Range("D97").Copy($ws.Range("E97"))
$ws.Cells.Item(97, 5).Value = 8981.85
$ws.Range("D98").Copy($ws.Range("E98"))
$ws.Range("D99").Copy($ws.Range("E99"))
$ws.Cells.Item(99, 5).Value = 56.94
$ws.Range("D100").Copy($ws.Range("E100"))
$ws.Cells.Item(100, 5).Value = 12232.61
$ws.Range("D101").Copy($ws.Range("E101"))
$ws.Range("D102").Copy($ws.Range("E102"))
$ws.Cells.Item(102, 5).Value = 5.94
$ws.Range("D103").Copy($ws.Range("E103"))
$ws.Cells.Item(103, 5).Value = 2274.27
$ws.Range("D104").Copy($ws.Range("E104"))
$ws.Range("D105").Copy($ws.Range("E105"))
$ws.Cells.Item(105, 5).Value = 26.15
$ws.Range("D106").Copy($ws.Range("E106"))
$ws.Cells.Item(106, 5).Value = 873.12
$ws.Range("D107").Copy($ws.Range("E107"))
$ws.Range("D108").Copy($ws.Range("E108"))
$ws.Cells.Item(108, 5).Value = 29.58
$ws.Range("D109").Copy($ws.Range("E109"))
$ws.Cells.Item(109, 5).Value = 2734.02
$ws.Range("D110").Copy($ws.Range("E110"))
$ws.Range("D111").Copy($ws.Range("E111"))
$ws.Cells.Item(111, 5).Value = 21.2
$ws.Range("D112").Copy($ws.Range("E112"))
$ws.Cells.Item(112, 5).Value = 3971.97
$ws.Range("D113").Copy($ws.Range("E113"))
$ws.Range("D114").Copy($ws.Range("E114"))
$ws.Cells.Item(114, 5).Value = 29.02
$ws.Range("D115").Copy($ws.Range("E115"))
$ws.Cells.Item(115, 5).Value = 3358.45
